$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = -13.027
$ws.Range("B10").Value = 5.599
$ws.Range("B12").Value = 5.315
$ws.Range("C15").Value = -13.895
$ws.Range("B18").Value = 5.154
$ws.Range("C20").Value = -12.46
$ws.Range("C29").Value = -12.1
$ws.Range("C30").Value = -12.952
$ws.Range("C31").Value = -13.024
$ws.Range("B37").Value = 8.494
$ws.Range("C40").Value = -12.782
$ws.Range("B55").Value = 5.315
$ws.Range("B68").Value = 5.431
$ws.Range("C68").Value = -11.177
$ws.Range("C76").Value = -12.977
$ws.Range("B77").Value = 6.377
$ws.Range("B78").Value = 7.631
$ws.Range("C87").Value = -13.288
$ws.Range("C88").Value = -13.35
$ws.Range("C96").Value = -12.57
$ws.Range("C98").Value = -13.23
$ws.Range("C101").Value = -12.612
$ws.Range("C102").Value = -13.091
